$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.448.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.26%  "
$ws.Range("D3").Value = "'3.256.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'554.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").Value = "'181.30"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "'3.255.17"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.82%  "
$ws.Range("E10").Value = "  -8.03%  "
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").Value = "'47.08"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.82%  "
$ws.Range("E13").Value = "  -6.12%  "
$ws.Range("D14").Value = "'631.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'8.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").Value = "'3.781.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").Value = "'65.393.09"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.13%  "
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "'3.252.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.50%  "
$ws.Range("D21").Value = "'11.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.39%  "
$ws.Range("D22").Value = "'0.899"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").Value = "'17.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "'105.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.78%  "
$ws.Range("D25").Value = "'4.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.07%  "
$ws.Range("D26").Value = "'3.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.99%  "
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("D28").Value = "'9.50"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "'8.65"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").Value = "'30.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.38%  "
$ws.Range("D31").Value = "'3.98"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("D32").Value = "'6.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("D33").Value = "'10.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").Value = "'543.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.01%  "
$ws.Range("D35").Value = "'0.104"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'56.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.62%  "
$ws.Range("D38").Value = "'3.594.71"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'3.38"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").Value = "'2.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").Value = "'0.0₃0712"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.75%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'3.25"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.84%  "
$ws.Range("D44").Value = "'31.82"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.333"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.39%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "'0.0413"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").Value = "'2.59"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.42%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("E51").Value = "  +1.66%  "
